$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; existing rows shift down by one.
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "地点"
$ws.Range("B1").Value = "开始时间"
$ws.Range("C1").Value = "结束时间"
$ws.Range("D1").Value = "准考证号"
$ws.Range("E1").Value = "姓名"
$ws.Range("F1").Value = "邮箱"
$ws.Range("G1").Value = "电话"

$ws.Range("H12").Select()

